# Update the two-digit division problems in the worksheet table.
$d = $word.ActiveDocument

$d.Content.Find.Execute("86÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷5=", 2) | Out-Null
$d.Content.Find.Execute("25÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷4=", 2) | Out-Null
$d.Content.Find.Execute("85÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷3=", 2) | Out-Null
$d.Content.Find.Execute("51÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷7=", 2) | Out-Null
$d.Content.Find.Execute("99÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷8=", 2) | Out-Null
$d.Content.Find.Execute("19÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "38÷8=", 2) | Out-Null
$d.Content.Find.Execute("70÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷3=", 2) | Out-Null
$d.Content.Find.Execute("80÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷2=", 2) | Out-Null
$d.Content.Find.Execute("49÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷6=", 2) | Out-Null
$d.Content.Find.Execute("51÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷6=", 2) | Out-Null
$d.Content.Find.Execute("69÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷8=", 2) | Out-Null
$d.Content.Find.Execute("61÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷5=", 2) | Out-Null
$d.Content.Find.Execute("57÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷8=", 2) | Out-Null
$d.Content.Find.Execute("45÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷3=", 2) | Out-Null
$d.Content.Find.Execute("54÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "77÷5=", 2) | Out-Null
$d.Content.Find.Execute("47÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷2=", 2) | Out-Null
$d.Content.Find.Execute("16÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷5=", 2) | Out-Null
$d.Content.Find.Execute("95÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷8=", 2) | Out-Null
$d.Content.Find.Execute("55÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷3=", 2) | Out-Null
$d.Content.Find.Execute("53÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "96÷7=", 2) | Out-Null
$d.Content.Find.Execute("21÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷6=", 2) | Out-Null
$d.Content.Find.Execute("22÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "94÷3=", 2) | Out-Null
$d.Content.Find.Execute("21÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷9=", 2) | Out-Null
$d.Content.Find.Execute("47÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷5=", 2) | Out-Null
$d.Content.Find.Execute("79÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷6=", 2) | Out-Null
